$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 95 in column C become 7586
$ws.Range("C2:C95").Value = 7586

# Rows 96 through 252 in column C become 7569
$ws.Range("C96:C252").Value = 7569
